$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($headerOrFooter, $newName) {
    $inlineShape = $headerOrFooter.Range.InlineShapes.Item(1)
    # Route the rename through Shape.Name (ConvertToShape / ConvertToInlineShape)
    # rather than the bare InlineShape, since that is the reliably-wired
    # property path for both header and footer stories.
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    return $shape.ConvertToInlineShape()
}

# Pearson logo (PearsonLogo.png) inline pictures sit in both footers
# (primary/default + first-page) and are renamed image1.png -> image2.png.
Rename-InlineLogo $sec.Footers.Item(1) "image2.png" | Out-Null
Rename-InlineLogo $sec.Footers.Item(2) "image2.png" | Out-Null

# BTec logo (BTec_Logo-Orange) inline pictures sit in both headers
# (primary/default + first-page) and are renamed image2.jpg -> image1.jpg.
Rename-InlineLogo $sec.Headers.Item(1) "image1.jpg" | Out-Null
Rename-InlineLogo $sec.Headers.Item(2) "image1.jpg" | Out-Null
